$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 938.9091
$ws.Range("I41").Value = 2261.5
$ws.Range("J41").Value = 183.14285
$ws.Range("K41").Value = 2261.5
$ws.Range("L41").Value = 183.14285
$ws.Range("M41").Value = -1821.5
$ws.Range("N41").Value = -1063.14285
$ws.Range("H48").Value = 5000
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 15000
$ws.Range("H53").Value = 113
$ws.Range("J53").Value = 114.625
$ws.Range("L53").Value = 114.625
$ws.Range("N53").Value = -1388.625
$ws.Range("H56").Value = 5000
$ws.Range("J56").Value = 5000
$ws.Range("L56").Value = 15000
$ws.Range("H86").Value = 4548.6665
$ws.Range("J86").Value = 7999
$ws.Range("L86").Value = 7999
$ws.Range("N86").Value = -10245
$ws.Range("H89").Value = 4548.6665
$ws.Range("J89").Value = 7999
$ws.Range("L89").Value = 39995
$ws.Range("N89").Value = -51227
$ws.Range("H98").Value = 3886.6
$ws.Range("I98").Value = 3273.8
$ws.Range("J98").Value = 4499.4
$ws.Range("K98").Value = 3273.8
$ws.Range("L98").Value = 4499.4
$ws.Range("M98").Value = -1775.8
$ws.Range("N98").Value = -7495.4
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H110").Value = 39999
$ws.Range("J110").Value = 39999
$ws.Range("L110").Value = 39999
$ws.Range("H113").Value = 2637.2
$ws.Range("J113").Value = 2499
$ws.Range("L113").Value = 2499
$ws.Range("N113").Value = -9007
$ws.Range("H122").Value = 3886.6
$ws.Range("I122").Value = 3273.8
$ws.Range("J122").Value = 4499.4
$ws.Range("K122").Value = 9821.400000000001
$ws.Range("L122").Value = 13498.2
$ws.Range("M122").Value = -7371.400000000001
$ws.Range("N122").Value = -18398.2
$ws.Range("H125").Value = 1119.8
$ws.Range("I125").Value = 1119.8
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 10078.2
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -7618.199999999999
$ws.Range("H129").Value = 2078.6667
$ws.Range("I129").Value = 2078.6667
$ws.Range("K129").Value = 6236.000100000001
$ws.Range("M129").Value = -1236.000100000001
$ws.Range("H135").Value = 6651.9546
$ws.Range("I135").Value = 1519.1111
$ws.Range("K135").Value = 13671.9999
$ws.Range("M135").Value = -11136.9999
$ws.Range("H138").Value = 2894.3489
$ws.Range("I138").Value = 1471.3077
$ws.Range("K138").Value = 4413.9231
$ws.Range("M138").Value = 726.0769
$ws.Range("N48").Value = -15584
$ws.Range("N56").Value = -16068
$ws.Range("N110").Value = -48179
$ws.Range("N125").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 26820.625
$ws.Range("I16").Value = 26820.625
$ws.Range("K16").Value = 26820.625
$ws.Range("M16").Value = -26533.625
$ws.Range("H32").Value = 2679.2292
$ws.Range("I32").Value = 1228.738
$ws.Range("K32").Value = 1228.738
$ws.Range("M32").Value = -941.7380000000001
$ws.Range("H113").Value = 19996.334
$ws.Range("J113").Value = 19996.334
$ws.Range("L113").Value = 19996.334
$ws.Range("N113").Value = -28674.334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2753.75
$ws.Range("I99").Value = 2753.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2753.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1255.75
$ws.Range("N99").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3428.537
$ws.Range("J31").Value = 2463.6316
$ws.Range("L31").Value = 2463.6316
$ws.Range("N31").Value = -3053.6316
$ws.Range("H34").Value = 3428.537
$ws.Range("J34").Value = 2463.6316
$ws.Range("L34").Value = 2463.6316
$ws.Range("N34").Value = -2867.6316
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("H99").Value = 10198.875
$ws.Range("I99").Value = 5699.3335
$ws.Range("K99").Value = 5699.3335
$ws.Range("M99").Value = -4201.3335
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H126").Value = 10198.875
$ws.Range("I126").Value = 5699.3335
$ws.Range("K126").Value = 17098.0005
$ws.Range("M126").Value = -14628.0005
$ws.Range("M76").ClearContents()
$ws.Range("M79").ClearContents()
$ws.Range("N122").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 2197.6667
$ws.Range("I44").Value = 437.4
$ws.Range("K44").Value = 1312.2
$ws.Range("M44").Value = -914.1999999999998
$ws.Range("H69").Value = 2198
$ws.Range("J69").Value = 1996
$ws.Range("L69").Value = 5988
$ws.Range("N69").Value = -7610
$ws.Range("H72").Value = 2198
$ws.Range("J72").Value = 1996
$ws.Range("L72").Value = 17964
$ws.Range("N72").Value = -26076
$ws.Range("H122").Value = 2351.6843
$ws.Range("J122").Value = 3177.5386
$ws.Range("L122").Value = 28597.8474
$ws.Range("N122").Value = -33497.8474
$ws.Range("H131").Value = 2576
$ws.Range("J131").Value = 2707.3333
$ws.Range("L131").Value = 8121.999899999999
$ws.Range("N131").Value = -18201.9999
$ws.Range("H141").Value = 3414
$ws.Range("I141").Value = 3414
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 10242
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -5062
$ws.Range("N141").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2799.8
$ws.Range("I102").Value = 2874.75
$ws.Range("K102").Value = 2874.75
$ws.Range("M102").Value = -1252.75
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H126").Value = 3283.625
$ws.Range("I126").Value = 2515.75
$ws.Range("K126").Value = 7547.25
$ws.Range("M126").Value = -5077.25
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 42220.6
$ws.Range("I7").Value = 42220.6
$ws.Range("K7").Value = 42220.6
$ws.Range("M7").Value = -42108.6
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("H122").Value = 7952.4165
$ws.Range("I122").Value = 10467.25
$ws.Range("J122").Value = 2922.75
$ws.Range("K122").Value = 31401.75
$ws.Range("L122").Value = 8768.25
$ws.Range("M122").Value = -28951.75
$ws.Range("N122").Value = -13668.25
$ws.Range("H126").Value = 42220.6
$ws.Range("I126").Value = 42220.6
$ws.Range("K126").Value = 126661.8
$ws.Range("M126").Value = -124191.8
$ws.Range("M40").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18000
$ws.Range("I62").Value = 18000
$ws.Range("K62").Value = 18000
$ws.Range("M62").Value = -17376
$ws.Range("H65").Value = 18000
$ws.Range("I65").Value = 18000
$ws.Range("K65").Value = 90000
$ws.Range("M65").Value = -86880
$ws.Range("H122").Value = 3030.5454
$ws.Range("I122").Value = 3116.1667
$ws.Range("K122").Value = 9348.500100000001
$ws.Range("M122").Value = -6898.500100000001
$ws.Range("H126").Value = 1999.3334
$ws.Range("I126").Value = 1999.3334
$ws.Range("K126").Value = 5998.0002
$ws.Range("M126").Value = -3528.0002
